$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "Outstanding" (heading) / "Outstanding" columns one place to
# the right (N->O, O->P, P->Q).
$ws.Columns("N").Insert()

# Give the newly inserted column the same width as column M (11 chars)
# but without Excel's "best fit" auto-size flag.
$ws.Columns("N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and select cell S6,
# matching the saved selection state.
$ws.Activate()
$ws.Range("S6").Select()
